$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("EV Home win")

# Remove the AUSTRALIA / QUEENSLAND PREMIER LEAGUE match row (row 4),
# shifting all subsequent rows up by one.
$ws.Rows.Item(4).Delete()
